$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Docente(s) Responsável(eis)") {
        $targetPara = $p
        break
    }
}

# Insert a new, empty paragraph right after the heading paragraph
$targetPara.Range.InsertParagraphAfter()

# Re-fetch the heading paragraph (object references can get stale after
# the document structure changes) and grab the paragraph that now follows it
$insertedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Docente(s) Responsável(eis)") {
        $insertedPara = $p.Next()
        break
    }
}

# Give the new paragraph the requested text and the "List Bullet" style
$insertedPara.Range.Text = "8855158 - Morun Bernardino Neto"
$insertedPara.Style = $d.Styles.Item("List Bullet")
